# Added test row for fm25 - "Direct level from fm24 testing % TIV deductibles"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ftests")
$ws.Activate()

# Copy the right-aligned number format used by the Levels/Layers columns
# (F28:G28) onto the new row's F30:G30 cells before writing the values,
# so the new cells pick up the existing "right" alignment style instead
# of creating a brand-new style entry.
$ws.Range("F28:G28").Copy()
$ws.Range("F30:G30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C30").Value = "Direct level from fm24 testing % TIV deductibles"
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = "in progress"
$ws.Range("I30").Value = "in progress"
